$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.217.95"
$ws.Range("E2").Value = "  -3.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.463.25"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.11"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.95"
$ws.Range("E6").Value = "  -6.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  -2.77%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -4.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.31"
$ws.Range("E10").Value = "  -6.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0777"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.97"
$ws.Range("E13").Value = "  -4.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.846.20"
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.488.92"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.81"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.172.10"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  -5.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  -7.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.28"
$ws.Range("E22").Value = "  -1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.40"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -5.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.00"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("E28").Value = "  -6.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.60"
$ws.Range("E29").Value = "  -5.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.44"
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.48"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.47"
$ws.Range("E32").Value = "  -5.24%  "
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.56"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0742"
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.04"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.93"
$ws.Range("E38").Value = "  -7.98%  "
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("E40").Value = "  -7.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.95"
$ws.Range("E43").Value = "  -10.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.973.15"
$ws.Range("E44").Value = "  -1.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  -5.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("E46").Value = "  -7.91%  "
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "69.13"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.93"
$ws.Range("E49").Value = "  -3.74%  "
$ws.Range("E50").Value = "  -6.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.95"
$ws.Range("E51").Value = "  -6.70%  "
